# 供水用水情况.xlsx — refresh the time series:
#   - drop the oldest 6 years (2004年-2009年), which were rows 2-7
#   - this shifts the existing 2010年-2020年 rows up to rows 2-12
#   - append 2021年 (full data) as the new row 13
#   - append 2022年 (only 人均用水量/B and 用水总量/K populated so far) as row 14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2004年-2009年 rows entirely; everything below shifts up 6 rows.
$ws.Rows("2:7").Delete()

# Prime rows 13 and 14 with the same look (font/border/alignment) as the
# rest of the year column before filling in their values.
$ws.Range("A2").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A14").PasteSpecial(-4122)

# Row 13: 2021年
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 419.170749667224
$ws.Range("C13").Value = 5920.2
$ws.Range("D13").Value = 138.3
$ws.Range("E13").Value = 3644.3
$ws.Range("F13").Value = 853.8
$ws.Range("G13").Value = 4928.1
$ws.Range("H13").Value = 1049.6
$ws.Range("I13").Value = 316.9
$ws.Range("J13").Value = 909.4
$ws.Range("K13").Value = 5920.2

# Row 14: 2022年 — only per-capita water use and total water use are in yet.
$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = 425
$ws.Range("K14").Value = 5997

$wb.Save()
